# Added filtering options for the Component Analysis
# Clears the trailing (higher-horizon) error columns that are no longer
# needed for the first evaluation, row by row, matching the staircase
# pattern already used by later rows in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ranges = @(
    "I2:K2","G3:K3","I4:K4","G5:K5","I6:K6","G7:K7","I8:K8","G9:K9",
    "I10:K10","G11:K11","I12:K12","K13:K13","I14:K14","K15:K15","I16:K16",
    "K17:K17","J18:K18","I19:K19","K21:K21","J22:K22","I23:K23","K25:K25",
    "J26:K26","I27:K27","K29:K29","J30:K30","I31:K31","K33:K33","J34:K34",
    "I35:K35","K37:K37","J38:K38","I39:J39"
)

foreach ($addr in $ranges) {
    $ws.Range($addr).ClearContents()
}
